$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new 2021 column (M) ------------------------------------------------

# Copy formatting from L2:L6 (the previous last column) into M2:M6 so the
# new column picks up borders/fonts/alignment consistently with its row
# before we put any values in it.
$ws.Range("L2:L6").Copy()
$ws.Range("M2:M6").PasteSpecial(-4122)  # xlPasteFormats

# Header for the new column
$ws.Range("M3").Value = 2021

# New data values for 2021
$ws.Range("M4").Value = 7105
$ws.Range("M5").Value = 81079
$ws.Range("M6").Value = 214139

# --- Number-format clean-up -------------------------------------------------
# The data columns (2012-2020, D:L) used a custom "0.0" number format; the
# refreshed sheet drops that custom format in favor of the default General
# format for all data cells, including the newly added 2021 column.
$ws.Range("D4:M6").NumberFormat = "general"

# --- Misc sheet-view clean-up ------------------------------------------------
# Reset the active cell / clear the stale selection that used to be saved
# with the sheet view.
$ws.Range("A1").Select()
